$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 12:05"

# Row 11 - Alemania
$ws.Range("B11").Value = 179110
$ws.Range("C11").Value = 89
$ws.Range("E11").Value = 11800
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 8310

# Row 40 - Rumania
$ws.Range("B40").Value = 17712
$ws.Range("C40").Value = 127
$ws.Range("D40").Value = 10777
$ws.Range("E40").Value = 5776

# Row 58 - Marruecos
$ws.Range("B58").Value = 7300
$ws.Range("C58").Value = 89
$ws.Range("D58").Value = 4347
$ws.Range("E58").Value = 2756
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 197

# Row 64 - Finlandia
$ws.Range("B64").Value = 6537
$ws.Range("C64").Value = 44
$ws.Range("E64").Value = 1431

# Row 109 - Albania
$ws.Range("B109").Value = 981
$ws.Range("C109").Value = 12
$ws.Range("D109").Value = 777
$ws.Range("E109").Value = 173
